$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2505
$ws1.Range("F5").Value = 1419
$ws1.Range("F6").Value = 1110
$ws1.Range("F7").Value = 314
$ws1.Range("F8").Value = 522
$ws1.Range("F11").Value = 104
$ws1.Range("F13").Value = 8741
$ws1.Range("F16").Value = 264
$ws1.Range("F18").Value = 176
$ws1.Range("F19").Value = 471
$ws1.Range("F20").Value = 598
$ws1.Range("F22").Value = 1160
$ws1.Range("F24").Value = 2042
$ws1.Range("F25").Value = 2100
$ws1.Range("F27").Value = 1784
$ws1.Range("F29").Value = 1919
$ws1.Range("F31").Value = 219
$ws1.Range("F32").Value = 57
$ws1.Range("F33").Value = 108
$ws1.Range("F34").Value = 192
$ws1.Range("F35").Value = 10
$ws1.Range("F36").Value = 310
$ws1.Range("F38").Value = 259
$ws1.Range("F39").Value = 441
$ws1.Range("F40").Value = 742
$ws1.Range("F42").Value = 267

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2505
$ws4.Range("F5").Value = 1419
$ws4.Range("F7").Value = 1110
$ws4.Range("F8").Value = 314
$ws4.Range("F9").Value = 522
$ws4.Range("F12").Value = 104
$ws4.Range("F14").Value = 8741
$ws4.Range("F18").Value = 264
$ws4.Range("F20").Value = 176
$ws4.Range("F21").Value = 471
$ws4.Range("F22").Value = 598
$ws4.Range("F24").Value = 1160
$ws4.Range("F26").Value = 2042
$ws4.Range("F27").Value = 2100
$ws4.Range("F29").Value = 1784
$ws4.Range("F31").Value = 1919
$ws4.Range("F33").Value = 219
$ws4.Range("F34").Value = 57
$ws4.Range("F35").Value = 108
$ws4.Range("F36").Value = 192
$ws4.Range("F37").Value = 10
$ws4.Range("F38").Value = 310
$ws4.Range("F40").Value = 259
$ws4.Range("F41").Value = 441
$ws4.Range("F46").Value = 742
$ws4.Range("F49").Value = 267
